# Updated symbol list refresh: push new coinranking price/volume snapshot
# values into the existing "Price" (D) / "Volume(1h)" (E) columns.
#
# These columns store plain text (inline strings) in the workbook, e.g.
# "246.11", "0.00002101" - not numbers. If we just assign a numeric-looking
# string to Range.Value, Excel auto-converts it to a real number (and can
# silently drop meaningful trailing/leading zeros, e.g. "0.1430" -> 0.143).
# To faithfully reproduce a text update we temporarily force the cell to
# Text format before writing the value, then restore the original ("Normal")
# style so no formatting side effects remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "246.14"
Set-TextValue "D3"  "22.02"
Set-TextValue "D4"  "5.430"
Set-TextValue "D6"  "3.381"
Set-TextValue "D7"  "6.350"
Set-TextValue "D8"  "0.8082"
Set-TextValue "D9"  "0.9660"
Set-TextValue "D10" "0.1430"
Set-TextValue "D11" "0.07464"
Set-TextValue "D12" "0.03347"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCXBestin24h"
Set-TextValue "D13" "0.03034"
Set-TextValue "D14" "4.159"
Set-TextValue "D15" "0.09400"
Set-TextValue "D16" "0.001589"
Set-TextValue "D17" "0.04821"
Set-TextValue "D18" "0.0005891"
Set-TextValue "D19" "0.006101"
Set-TextValue "D20" "0.004107"
Set-TextValue "D21" "0.0009954"
Set-TextValue "D24" "2.212"
Set-TextValue "D25" "0.3210"
Set-TextValue "D40" "0.03865"
Set-TextValue "D41" "0.006711"
Set-TextValue "D42" "0.1078"
Set-TextValue "D43" "0.002540"
Set-TextValue "E43" "42CEJICEJI"
Set-TextValue "D44" "0.006688"
Set-TextValue "D48" "0.1463"
Set-TextValue "D49" "0.00002100"
